# Fix newline formatting in header cells: replace the literal two-character
# sequence "\n" (backslash + n) that was baked into the text with an actual
# line-break character (Chr(10)), so wrapped text renders as two lines.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fixes = @{
    "G2" = "B3:0/0`nFire Alarm Zone 1"
    "H2" = "B3:0/8`nFire Eye Failure Warning"
    "I2" = "B3:0/9`nStrobe Light On"
    "J2" = "B3:0/10`nStrobe Light On"
    "K2" = "B3:0/11`nFire Alarm Zone 2"
    "L2" = "B3:0/13`nStrobe Light On"
    "M2" = "B3:2/0`nPlant ESD"
    "N2" = "B3:2/4`nPlant ESD"
    "O2" = "B3:2/8`nPlant ESD"
    "P2" = "B3:2/9`nPlant ESD"
    "Q2" = "B3:2/10`nPlant ESD"
    "R2" = "B3:2/13`nPlant ESD"
    "S2" = "B3:10/0`nESD Alarm to Office PLC"
    "T2" = "O:0/0`nDeluge Valve Zone 2 Open"
    "A3" = "Interlock`nNo"
    "E3" = "Pre-Trip`n(H or L)"
    "F3" = "Trip`n(HH or LL)"
}

foreach ($addr in $fixes.Keys) {
    $ws.Range($addr).Value = $fixes[$addr]
}
